$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet
$ws.Name = "Through 2022-12-15"

# Update header label cell (I1) which references the "2022 (through 12-14)" shared string
$ws.Range("I1").Value = "2022 (through 12-15)"

# Update December total (I13) and yearly total (I14)
$ws.Range("I13").Value = 64
$ws.Range("I14").Value = 1580
